# Insert a new data row at row 486 (pushing the existing rows 486:571 down
# to 487:572) and populate it with a new Betarraga / Vega Modelo de Temuco
# price record, per the weekly Fruta/Hortaliza update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(486).Insert()

$ws.Cells.Item(486, 1).Value  = 10
$ws.Cells.Item(486, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(486, 3).Value  = "La Araucanía"
$ws.Cells.Item(486, 4).Value  = 45034
$ws.Cells.Item(486, 5).Value  = 9
$ws.Cells.Item(486, 6).Value  = 100114014
$ws.Cells.Item(486, 7).Value  = "Betarraga"
$ws.Cells.Item(486, 8).Value  = "Sin especificar"
$ws.Cells.Item(486, 9).Value  = "Primera"
$ws.Cells.Item(486, 10).Value = 65
$ws.Cells.Item(486, 11).Value = 12000
$ws.Cells.Item(486, 12).Value = 12000
$ws.Cells.Item(486, 13).Value = 12000
$ws.Cells.Item(486, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(486, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(486, 16).Value = 1000
$ws.Cells.Item(486, 17).Value = 12
$ws.Cells.Item(486, 18).Value = "Hortaliza"
